$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update "Valor Mora" header total (E11) ---
$ws.Range("E11").Value = 991413

# --- 2. Update "Cant. Periodos" (F13): 21 -> 22 ---
$ws.Range("F13").Value = 22

# --- 3. Insert a new row at 37 for the new period 2508. This shifts the
#        trailing rows (old 41/42 signature block) down to 42/43. ---
$ws.Rows("37").Insert()

# Duplicate the formatting (and values, overwritten below) of the previous
# "last" row (36) into the newly inserted row 37.
$ws.Range("B36:J36").Copy($ws.Range("B37:J37"))

# --- 4. Re-apply "regular" data-row formatting (taken from row 35) to the
#        row that used to be the last one (row 36), since it is no longer
#        the last period row. ---
$ws.Range("B35:J35").Copy($ws.Range("B36:J36"))

# --- 5. Re-write the period labels (column E) in ascending order, and the
#        "Valor Mora" amounts (column F) -- all 46400 except period 2311,
#        which carries 17013. ---
$periods = @("2311","2312","2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412","2501","2502","2503","2504","2505","2506","2507","2508")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $p = $periods[$i]
    $ws.Range("E$r").Value = $p
    if ($p -eq "2311") {
        $ws.Range("F$r").Value = 17013
    } else {
        $ws.Range("F$r").Value = 46400
    }
}

Write-Host "Done"
